$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 20.79417744437257
$ws.Range("D2").Value = 9.04138735040676
$ws.Range("E2").Value = 11.8011012734251
$ws.Range("F2").Value = 57.43733929506016
$ws.Range("G2").Value = 3.787086984313411
$ws.Range("I2").Value = 31.55432770389527
$ws.Range("J2").Value = 8.196796425028694
$ws.Range("L2").Value = 17.98559190041663
$ws.Range("M2").Value = 22.60695363921343
$ws.Range("N2").Value = 20.91452401751551
$ws.Range("B3").Value = 20.57391262408368
$ws.Range("D3").Value = 8.719123668152823
$ws.Range("E3").Value = 11.22245571595289
$ws.Range("F3").Value = 57.23128565453276
$ws.Range("G3").Value = 3.793828923513315
$ws.Range("I3").Value = 31.68962779581675
$ws.Range("J3").Value = 7.953951839490816
$ws.Range("L3").Value = 17.85544093948653
$ws.Range("M3").Value = 22.49046686374407
$ws.Range("N3").Value = 20.99954715728993
$ws.Range("B4").Value = 20.44263337141837
$ws.Range("D4").Value = 8.518982885716426
$ws.Range("E4").Value = 10.85579203637742
$ws.Range("F4").Value = 57.12309779356739
$ws.Range("G4").Value = 3.798174478152341
$ws.Range("I4").Value = 31.77917511379043
$ws.Range("J4").Value = 7.799718502083307
$ws.Range("L4").Value = 17.77996660876455
$ws.Range("M4").Value = 22.42402495755564
$ws.Range("N4").Value = 21.05382993037099
$ws.Range("B5").Value = 20.39018679306658
$ws.Range("D5").Value = 8.43699609545247
$ws.Range("E5").Value = 10.70373499692437
$ws.Range("F5").Value = 57.08362414028504
$ws.Range("G5").Value = 3.799997385190803
$ws.Range("I5").Value = 31.81728811360954
$ws.Range("J5").Value = 7.735624617385644
$ws.Range("L5").Value = 17.75034846591914
$ws.Range("M5").Value = 22.3982442382193
$ws.Range("N5").Value = 21.07647607422461
$ws.Range("B6").Value = 20.38154307721462
$ws.Range("D6").Value = 8.423360800648783
$ws.Range("E6").Value = 10.67833379375273
$ws.Range("F6").Value = 57.07734823920845
$ws.Range("G6").Value = 3.800303229422412
$ws.Range("I6").Value = 31.82371452996261
$ws.Range("J6").Value = 7.724908267791773
$ws.Range("L6").Value = 17.74549974719656
$ws.Range("M6").Value = 22.39404204186771
$ws.Range("N6").Value = 21.0802682707867
$ws.Range("B7").Value = 20.44192173424102
$ws.Range("D7").Value = 8.517878718522946
$ws.Range("E7").Value = 10.85375170892853
$ws.Range("F7").Value = 57.12254675564184
$ws.Range("G7").Value = 3.798198851381196
$ws.Range("I7").Value = 31.77968256019806
$ws.Range("J7").Value = 7.798859076227499
$ws.Range("L7").Value = 17.77956253174774
$ws.Range("M7").Value = 22.42367200549284
$ws.Range("N7").Value = 21.05413321273123
$ws.Range("B8").Value = 20.71744034986412
$ws.Range("D8").Value = 8.930824214985423
$ws.Range("E8").Value = 11.60406744567581
$ws.Range("F8").Value = 57.36249227378616
$ws.Range("G8").Value = 3.789369013856358
$ws.Range("I8").Value = 31.59963226265539
$ws.Range("J8").Value = 8.114153973744019
$ws.Range("L8").Value = 17.93980842512169
$ws.Range("M8").Value = 22.56574517120109
$ws.Range("N8").Value = 20.94341041247094
$ws.Range("B9").Value = 21.28672384047691
$ws.Range("D9").Value = 9.716895678532142
$ws.Range("E9").Value = 12.97706305062472
$ws.Range("F9").Value = 57.97797858731014
$ws.Range("G9").Value = 3.773676008539018
$ws.Range("I9").Value = 31.29817749112257
$ws.Range("J9").Value = 8.68995590254276
$ws.Range("L9").Value = 18.28820666869341
$ws.Range("M9").Value = 22.88383499013241
$ws.Range("N9").Value = 20.74263892744421
$ws.Range("B10").Value = 21.71925506777001
$ws.Range("D10").Value = 10.27315504747319
$ws.Range("E10").Value = 14.02992524221436
$ws.Range("F10").Value = 58.51749130059397
$ws.Range("G10").Value = 3.763118428879909
$ws.Range("I10").Value = 31.10854512051655
$ws.Range("J10").Value = 9.084952400392872
$ws.Range("L10").Value = 18.56345067845115
$ws.Range("M10").Value = 23.14037908281028
$ws.Range("N10").Value = 20.60491152080493
$ws.Range("B11").Value = 21.91838603907587
$ws.Range("D11").Value = 10.52038842196682
$ws.Range("E11").Value = 14.51174450410116
$ws.Range("F11").Value = 58.78153726247359
$ws.Range("G11").Value = 3.758522934662111
$ws.Range("I11").Value = 31.02927617707964
$ws.Range("J11").Value = 9.258181797125069
$ws.Range("L11").Value = 18.69249032627149
$ws.Range("M11").Value = 23.26175448996871
$ws.Range("N11").Value = 20.54433760217839
$ws.Range("B12").Value = 21.99407032911045
$ws.Range("D12").Value = 10.61308814642016
$ws.Range("E12").Value = 14.69019916671644
$ws.Range("F12").Value = 58.88416279403265
$ws.Range("G12").Value = 3.756812248150441
$ws.Range("I12").Value = 31.0002724779213
$ws.Range("J12").Value = 9.322824079717291
$ws.Range("L12").Value = 18.74187267629946
$ws.Range("M12").Value = 23.30836118953837
$ws.Range("N12").Value = 20.52169555084698
$ws.Range("B13").Value = 21.97775901242171
$ws.Range("D13").Value = 10.59316588219547
$ws.Range("E13").Value = 14.65194331616295
$ws.Range("F13").Value = 58.86194396546055
$ws.Range("G13").Value = 3.75717936613768
$ws.Range("I13").Value = 31.00647373188242
$ws.Range("J13").Value = 9.308945125412585
$ws.Range("L13").Value = 18.73121480587919
$ws.Range("M13").Value = 23.29829538810541
$ws.Range("N13").Value = 20.52655880565471
$ws.Range("B14").Value = 21.92460742558023
$ws.Range("D14").Value = 10.52803382404183
$ws.Range("E14").Value = 14.52650611717183
$ws.Range("F14").Value = 58.78992772656756
$ws.Range("G14").Value = 3.758381605220463
$ws.Range("I14").Value = 31.02686966355737
$ws.Range("J14").Value = 9.26351926053766
$ws.Range("L14").Value = 18.69654285532116
$ws.Range("M14").Value = 23.26557609605541
$ws.Range("N14").Value = 20.54246891230944
$ws.Range("B15").Value = 21.89208484320707
$ws.Range("D15").Value = 10.48801607515597
$ws.Range("E15").Value = 14.44915194564705
$ws.Range("F15").Value = 58.74615783325217
$ws.Range("G15").Value = 3.759121848735167
$ws.Range("I15").Value = 31.03949501153059
$ws.Range("J15").Value = 9.235569357240173
$ws.Range("L15").Value = 18.67537173034231
$ws.Range("M15").Value = 23.24561768109305
$ws.Range("N15").Value = 20.55225276822145
$ws.Range("B16").Value = 21.70628368267207
$ws.Range("D16").Value = 10.25687318399132
$ws.Range("E16").Value = 13.9978763501128
$ws.Range("F16").Value = 58.50060661546893
$ws.Range("G16").Value = 3.763422901788345
$ws.Range("I16").Value = 31.1138669625275
$ws.Range("J16").Value = 9.073499103764622
$ws.Range("L16").Value = 18.55509190992372
$ws.Range("M16").Value = 23.13253861647324
$ws.Range("N16").Value = 20.60891174911508
$ws.Range("B17").Value = 21.59286277296454
$ws.Range("D17").Value = 10.1135210632649
$ws.Range("E17").Value = 13.71388667129672
$ws.Range("F17").Value = 58.3547107252278
$ws.Range("G17").Value = 3.766114337625082
$ws.Range("I17").Value = 31.16128897505657
$ws.Range("J17").Value = 8.972399953936026
$ws.Range("L17").Value = 18.48226195045063
$ws.Range("M17").Value = 23.06434582577894
$ws.Range("N17").Value = 20.64420056199183
$ws.Range("B18").Value = 21.52785294435244
$ws.Range("D18").Value = 10.03052615274379
$ws.Range("E18").Value = 13.54791029570702
$ws.Range("F18").Value = 58.27255054657755
$ws.Range("G18").Value = 3.767681899245607
$ws.Range("I18").Value = 31.18922286172149
$ws.Range("J18").Value = 8.913644103525845
$ws.Range("L18").Value = 18.4407345970444
$ws.Range("M18").Value = 23.02556462612632
$ws.Range("N18").Value = 20.66469361564045
$ws.Range("B19").Value = 21.50588262900336
$ws.Range("D19").Value = 10.00233512526097
$ws.Range("E19").Value = 13.4912607053232
$ws.Range("F19").Value = 58.2450351299203
$ws.Range("G19").Value = 3.768216009204129
$ws.Range("I19").Value = 31.19879357499934
$ws.Range("J19").Value = 8.893647143812748
$ws.Range("L19").Value = 18.42673739644045
$ws.Range("M19").Value = 23.01251061978739
$ws.Range("N19").Value = 20.67166594117684
$ws.Range("B20").Value = 21.60491359987276
$ws.Range("D20").Value = 10.12883795421058
$ws.Range("E20").Value = 13.74439021771613
$ws.Range("F20").Value = 58.37006020180062
$ws.Range("G20").Value = 3.765825811490137
$ws.Range("I20").Value = 31.15617266503249
$ws.Range("J20").Value = 8.983225089869954
$ws.Range("L20").Value = 18.48997754042393
$ws.Range("M20").Value = 23.07155956698792
$ws.Range("N20").Value = 20.64042375816089
$ws.Range("B21").Value = 21.94021229911135
$ws.Range("D21").Value = 10.54719035333031
$ws.Range("E21").Value = 14.56345844774591
$ws.Range("F21").Value = 58.8110094209998
$ws.Range("G21").Value = 3.758027679235457
$ws.Range("I21").Value = 31.020851308893
$ws.Range("J21").Value = 9.276888084517788
$ws.Range("L21").Value = 18.70671306030102
$ws.Range("M21").Value = 23.27516927850865
$ws.Range("N21").Value = 20.53778771941147
$ws.Range("B22").Value = 22.16093919352265
$ws.Range("D22").Value = 10.81519976441923
$ws.Range("E22").Value = 15.07546300382693
$ws.Range("F22").Value = 59.11454353816741
$ws.Range("G22").Value = 3.75310313200305
$ws.Range("I22").Value = 30.93832351467926
$ws.Range("J22").Value = 9.463231481019932
$ws.Range("L22").Value = 18.85136474949556
$ws.Range("M22").Value = 23.41198297195388
$ws.Range("N22").Value = 20.47243291731922
$ws.Range("B23").Value = 22.04300836130045
$ws.Range("D23").Value = 10.67267882009528
$ws.Range("E23").Value = 14.80432176470846
$ws.Range("F23").Value = 58.95115200635467
$ws.Range("G23").Value = 3.755715807976117
$ws.Range("I23").Value = 30.98182653998311
$ws.Range("J23").Value = 9.364295522732862
$ws.Range("L23").Value = 18.77389777906567
$ws.Range("M23").Value = 23.33862978032901
$ws.Range("N23").Value = 20.50715726834874
$ws.Range("B24").Value = 21.59946480149246
$ws.Range("D24").Value = 10.12191498729602
$ws.Range("E24").Value = 13.73060799741336
$ws.Range("F24").Value = 58.36311535140479
$ws.Range("G24").Value = 3.76595619115709
$ws.Range("I24").Value = 31.15848366141931
$ws.Range("J24").Value = 8.978333014207465
$ws.Range("L24").Value = 18.48648824960017
$ws.Range("M24").Value = 23.06829691265557
$ws.Range("N24").Value = 20.64213061227856
$ws.Range("B25").Value = 21.1299602606051
$ws.Range("D25").Value = 9.507530010213609
$ws.Range("E25").Value = 12.61718298662012
$ws.Range("F25").Value = 57.79603886122604
$ws.Range("G25").Value = 3.777749450831784
$ws.Range("I25").Value = 31.3741643063242
$ws.Range("J25").Value = 8.538971513620213
$ws.Range("L25").Value = 18.19045214965039
$ws.Range("M25").Value = 22.79367647109045
$ws.Range("N25").Value = 20.79522167326294
